$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The refreshed place list has one fewer row (60 data rows instead of 61):
#   - removed: CCA Gallery, CityWay Gallery, Ripple Mobile Arts
#   - added:   Indiana Arts Building, Lebanon Community Arts Council-N
# Drop the now-unused last row so dimension becomes A1:E60 ...
$ws.Rows("61:61").Delete()

# ... then rewrite every data row (2-60) with the refreshed values
# (business_status, name, rating, user_ratings_total, and the row's id in column A).
$data = @(
    @(41, 'OPERATIONAL', '10 East Arts HUB', 5, 8),
    @(15, 'OPERATIONAL', 'Acceleration Art and Photography', 4.6, 16),
    @(4, 'OPERATIONAL', 'Art On Main Gallery And Gifts', 5, 11),
    @(49, 'OPERATIONAL', 'Art To Remember', 3.9, 7),
    @(14, 'OPERATIONAL', 'ArtMix', 5, 8),
    @(0, 'OPERATIONAL', 'Arts Council of Indianapolis', 4.8, 12),
    @(5, 'OPERATIONAL', 'Arts for Lawrence', 4.9, 19),
    @(9, 'OPERATIONAL', 'Arts for Learning', 0, 0),
    @(51, 'OPERATIONAL', 'Carmel Art Education Studio', 5, 12),
    @(30, 'OPERATIONAL', 'Carmel International Arts Festival', 4.8, 18),
    @(19, 'OPERATIONAL', 'Cat Head Press: Printshop and Artist Cooperative', 5, 26),
    @(32, 'OPERATIONAL', 'Clowes Memorial Hall', 4.7, 797),
    @(40, 'OPERATIONAL', 'Corporate Fine Art Group', 5, 2),
    @(47, 'OPERATIONAL', 'Creative Energy Arts', 0, 0),
    @(12, 'OPERATIONAL', 'Creative Expressions Arts', 0, 0),
    @(39, 'OPERATIONAL', 'Cultural Arts Gallery', 5, 1),
    @(29, 'OPERATIONAL', 'Dance Arts', 4.8, 10),
    @(50, 'OPERATIONAL', 'Dance Magic Performing Arts', 0, 0),
    @(26, 'OPERATIONAL', 'Edison School Of The Arts', 3.6, 36),
    @(6, 'OPERATIONAL', 'Gallery 924', 4.9, 15),
    @(56, 'OPERATIONAL', 'Glass Arts Indiana, Inc.', 4.8, 16),
    @(20, 'OPERATIONAL', 'Harrison Center', 4.8, 134),
    @(33, 'OPERATIONAL', 'Herron School of Art and Design', 4.8, 20),
    @(21, 'OPERATIONAL', 'Herron School of Art and Design - Eskenazi Fine Arts Center', 4.5, 2),
    @(10, 'OPERATIONAL', 'High Frequency Arts', 5, 3),
    @(45, 'OPERATIONAL', 'Ignition Arts, LLC', 4.9, 8),
    @(58, 'OPERATIONAL', 'Indiana Arts Building', 4, 3),
    @(2, 'OPERATIONAL', 'Indiana Arts Commission', 0, 0),
    @(34, 'OPERATIONAL', 'Indiana Fine Arts Academy', 4, 4),
    @(38, 'OPERATIONAL', 'Indiana Performing Arts Centre', 4.8, 6),
    @(3, 'OPERATIONAL', 'Indianapolis Art Center', 4.7, 112),
    @(54, 'OPERATIONAL', 'Indianapolis Arts Chorale', 0, 0),
    @(52, 'OPERATIONAL', 'Indianapolis Arts Services', 5, 1),
    @(16, 'OPERATIONAL', 'Indianapolis Artsgarden', 4.6, 146),
    @(17, 'OPERATIONAL', 'Korka International Arts', 0, 0),
    @(42, 'OPERATIONAL', 'Landmark Keystone Art Cinema', 4.4, 513),
    @(59, 'OPERATIONAL', 'Lebanon Community Arts Council-N', 0, 0),
    @(13, 'OPERATIONAL', 'Magdalena Gallery of Arts', 5, 1),
    @(23, 'OPERATIONAL', 'Newfields', 4.7, 3896),
    @(11, 'OPERATIONAL', 'Nickel Plate Arts', 4.9, 41),
    @(31, 'OPERATIONAL', 'Palette Art Studio', 5, 7),
    @(46, 'OPERATIONAL', 'Performing Arts', 0, 0),
    @(36, 'OPERATIONAL', 'Phiri Art', 4.5, 4),
    @(53, 'OPERATIONAL', 'Red Barn Arts Collective', 5, 1),
    @(18, 'OPERATIONAL', 'Schrott Center for the Arts', 4.8, 106),
    @(1, 'OPERATIONAL', 'Sho Arts', 5, 1),
    @(57, 'OPERATIONAL', 'Sidney & Lois Eskenazi Fine Arts Center', 5, 1),
    @(35, 'OPERATIONAL', 'Steve Haigh Fine Art', 5, 1),
    @(27, 'OPERATIONAL', 'Studio Alchemy LLC', 5, 5),
    @(24, 'OPERATIONAL', 'The Art Studio of Carmel', 0, 0),
    @(7, 'OPERATIONAL', 'The Carmel Arts Council Children’s Art Gallery', 5, 1),
    @(22, 'OPERATIONAL', 'The Center for the Performing Arts', 4.7, 1193),
    @(8, 'OPERATIONAL', 'The Murphy Art Center', 5, 1),
    @(48, 'OPERATIONAL', 'The Palladium at the Center for the Performing Arts', 4.9, 197),
    @(28, 'OPERATIONAL', 'United Art & Education', 4.5, 132),
    @(25, 'OPERATIONAL', 'Warren Performing Arts Center', 4.5, 374),
    @(44, 'OPERATIONAL', 'Zionsville Performing Arts Center', 4.7, 175),
    @(37, 'OPERATIONAL', 'carmel academy of the arts', 0, 0),
    @(55, 'OPERATIONAL', 'iDREAM Music And Arts Center', 0, 0)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $r = $i + 2
    $row = $data[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $ws.Cells.Item($r, 5).Value = $row[4]
}